$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.157.12"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "2.160.64"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'235.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "'0.606"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.04%  "
$ws.Range("D7").Value = "'69.11"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.91%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "'0.565"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.09%  "
$ws.Range("D10").Value = "'38.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.19%  "
$ws.Range("D11").Value = "'0.0911"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.39%  "
$ws.Range("D12").Value = "'54.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.50%  "
$ws.Range("D13").Value = "'0.100"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("D14").Value = "'6.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.43%  "
$ws.Range("D15").Value = "2.493.45"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "'14.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "2.163.48"
$ws.Range("E17").Value = "  -2.31%  "
$ws.Range("D18").Value = "'0.782"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.68%  "
$ws.Range("D19").Value = "40.950.03"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").Value = "0.0₃0991"
$ws.Range("E20").Value = "  -7.58%  "
$ws.Range("D21").Value = "'70.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.24%  "
$ws.Range("D22").Value = "'5.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.32%  "
$ws.Range("D23").Value = "'223.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").Value = "'9.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -13.58%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").Value = "'1.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.87%  "
$ws.Range("D27").Value = "'10.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.70%  "
$ws.Range("E28").Value = "  -3.90%  "
$ws.Range("D29").Value = "'2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.92%  "
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").Value = "'167.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "'19.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.48%  "
$ws.Range("D33").Value = "'29.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("D34").Value = "'0.0751"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.49%  "
$ws.Range("D35").Value = "'5.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -10.84%  "
$ws.Range("D36").Value = "'0.119"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.30%  "
$ws.Range("D37").Value = "'0.101"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.64%  "
$ws.Range("D38").Value = "'4.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.20%  "
$ws.Range("D39").Value = "'0.0278"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.27%  "
$ws.Range("D40").Value = "'2.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("D41").Value = "'11.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -15.03%  "
$ws.Range("D42").Value = "'5.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.34%  "
$ws.Range("D43").Value = "'58.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.12%  "
$ws.Range("D44").Value = "'0.187"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.82%  "
$ws.Range("D45").Value = "'8.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.60%  "
$ws.Range("D46").Value = "'0.0955"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.76%  "
$ws.Range("D47").Value = "'96.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.34%  "
$ws.Range("D48").Value = "'1.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.72%  "
$ws.Range("D49").Value = "'1.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.87%  "
$ws.Range("D50").Value = "'2.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.57%  "
$ws.Range("E51").Value = "  -3.28%  "
